# Workbook "Fruta, Vega Modelo de Temuco - Papaya.xlsx"
# A new weekly record was inserted as row 65 (Especial quality, 2022-02-23),
# pushing the former rows 65-67 down to rows 66-68.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the current row 65, shifting rows 65-67 down to 66-68.
$ws.Rows.Item(65).Insert()

# Populate the newly inserted row 65 with the new weekly entry.
$ws.Range("A65").Value = 10
$ws.Range("B65").Value = "Vega Modelo de Temuco"
$ws.Range("C65").Value = "La Araucanía"
$ws.Range("D65").Value = 44615
$ws.Range("E65").Value = 9
$ws.Range("F65").Value = "Fruta"
$ws.Range("G65").Value = 100108
$ws.Range("H65").Value = "Tropicales y subtropicales"
$ws.Range("I65").Value = 100108004
$ws.Range("J65").Value = "Papaya"
$ws.Range("K65").Value = "Cultivar IV Región"
$ws.Range("L65").Value = "Especial"
$ws.Range("M65").Value = 80
$ws.Range("N65").Value = 25000
$ws.Range("O65").Value = 25000
$ws.Range("P65").Value = 25000
$ws.Range("Q65").Value = "$/bandeja 10 kilos"
$ws.Range("R65").Value = "Provincia del Elquí"
$ws.Range("S65").Value = 2500
$ws.Range("T65").Value = 10
